$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing single test-run row (row 2) method name from the
# placeholder "test" to the real test method that was executed.
$ws.Range("A2").Value = "verifyCustomerNavigationToRegistrationPage"

# Append the log rows for the other test methods that ran in this suite,
# each with the same "Passed" / "29-10-2024" execution data as row 2, and
# copy the green "Test Status" highlight fill from B2 onto each new
# B-column cell.
$methods = @(
    "verifyNewCustomerRegistrationSubmissionFlow",
    "verifyCustomerRegistrationAndLoginNavigation",
    "verifyCustomerEmailActivation",
    "verifyCustomerSuccessfulLogin"
)

$statusFill = $ws.Range("B2").Interior.Color

$row = 3
foreach ($method in $methods) {
    $ws.Cells.Item($row, 1).Value = $method
    $ws.Cells.Item($row, 2).Value = "Passed"
    $ws.Cells.Item($row, 2).Interior.Color = $statusFill
    $ws.Cells.Item($row, 3).Value = "29-10-2024"
    $row = $row + 1
}

# Widen column A so the longer method names fit (mirrors the bestFit
# recalculation Excel performs after the log save).
$ws.Columns.Item(1).ColumnWidth = 43.2
